$d = $word.ActiveDocument
Write-Output ("Paragraphs.Count: " + $d.Paragraphs.Count)
$p1 = $d.Paragraphs.Item(1)
Write-Output ("p1 range: " + $p1.Range.Start + "," + $p1.Range.End + " text=[" + $p1.Range.Text + "]")
